$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 3347.8

$ws.Range("H85").Value = 3347.8

$ws.Range("H129").Value = 2977280.5
$ws.Range("J129").Value = 1104.0264
$ws.Range("L129").Value = 3312.0792
$ws.Range("N129").Value = -13312.0792

$ws.Range("H132").Value = 5558621
$ws.Range("I132").Value = 6669489.5
$ws.Range("J132").Value = 4278.3335
$ws.Range("K132").Value = 20008468.5
$ws.Range("L132").Value = 12835.0005
$ws.Range("M132").Value = -20005938.5
$ws.Range("N132").Value = -17895.0005

$ws.Range("H135").Value = 1426.5714
$ws.Range("I135").Value = 1490
$ws.Range("K135").Value = 13410
$ws.Range("M135").Value = -10875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8498.646000000001
$ws.Range("I32").Value = 6012.7715
$ws.Range("J32").Value = 22999.584
$ws.Range("K32").Value = 6012.7715
$ws.Range("L32").Value = 22999.584
$ws.Range("M32").Value = -5725.7715
$ws.Range("N32").Value = -23573.584

$ws.Range("H61").Value = 3093.7222
$ws.Range("I61").Value = 1977.6428
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 1977.6428
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -1765.6428
$ws.Range("N61").Value = -7424

$ws.Range("H74").Value = 1812.1818
$ws.Range("I74").Value = 2255.6667
$ws.Range("K74").Value = 2255.6667
$ws.Range("M74").Value = -1381.6667

$ws.Range("H77").Value = 1812.1818
$ws.Range("I77").Value = 2255.6667
$ws.Range("K77").Value = 11278.3335
$ws.Range("M77").Value = -6910.333500000001

$ws.Range("H136").Value = 3093.7222
$ws.Range("I136").Value = 1977.6428
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 5932.928400000001
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -3382.928400000001
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2557.5
$ws.Range("I134").Value = 1422.7333
$ws.Range("K134").Value = 4268.199900000001
$ws.Range("M134").Value = -1733.199900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2354.6428
$ws.Range("I99").Value = 1736.875
$ws.Range("K99").Value = 1736.875
$ws.Range("M99").Value = -238.875

$ws.Range("H126").Value = 2354.6428
$ws.Range("I126").Value = 1736.875
$ws.Range("K126").Value = 5210.625
$ws.Range("M126").Value = -2740.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3181.5293
$ws.Range("I3").Value = 1747.2727
$ws.Range("K3").Value = 5241.8181
$ws.Range("M3").Value = -5129.8181

$ws.Range("H87").Value = 10618.272
$ws.Range("I87").Value = 7543
$ws.Range("K87").Value = 22629
$ws.Range("M87").Value = -21381

$ws.Range("H90").Value = 10618.272
$ws.Range("I90").Value = 7543
$ws.Range("K90").Value = 67887
$ws.Range("M90").Value = -61647

$ws.Range("H102").Value = 1633.3334
$ws.Range("J102").Value = 2300
$ws.Range("L102").Value = 6900
$ws.Range("N102").Value = -11768

$ws.Range("H128").Value = 174666.67
$ws.Range("I128").Value = 174666.67
$ws.Range("K128").Value = 524000.01
$ws.Range("M128").Value = -519020.01

$ws.Range("H129").Value = 28934.316
$ws.Range("I129").Value = 4118.5713
$ws.Range("J129").Value = 43410.168
$ws.Range("K129").Value = 12355.7139
$ws.Range("L129").Value = 130230.504
$ws.Range("M129").Value = -7355.713899999999
$ws.Range("N129").Value = -140230.504

$ws.Range("H133").Value = 5253.8335
$ws.Range("J133").Value = 3926
$ws.Range("L133").Value = 11778
$ws.Range("N133").Value = -21898

$ws.Range("H134").Value = 2471.8125
$ws.Range("J134").Value = 3966.5
$ws.Range("L134").Value = 11899.5
$ws.Range("N134").Value = -22039.5

$ws.Range("H137").Value = 2968.16
$ws.Range("J137").Value = 3564.923
$ws.Range("L137").Value = 10694.769
$ws.Range("N137").Value = -20894.769

$ws.Range("H138").Value = 2798.8572
$ws.Range("J138").Value = 5444
$ws.Range("L138").Value = 16332
$ws.Range("N138").Value = -26612

$ws.Range("H139").Value = 3568.1765
$ws.Range("I139").Value = 3350.9092
$ws.Range("J139").Value = 3966.5
$ws.Range("K139").Value = 10052.7276
$ws.Range("L139").Value = 11899.5
$ws.Range("M139").Value = -4912.7276
$ws.Range("N139").Value = -22179.5

$ws.Range("H140").Value = 1983.5714
$ws.Range("J140").Value = 4760
$ws.Range("L140").Value = 14280
$ws.Range("N140").Value = -24640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1928.75
$ws.Range("I97").Value = 720
$ws.Range("J97").Value = 3137.5
$ws.Range("K97").Value = 720
$ws.Range("L97").Value = 3137.5
$ws.Range("M97").Value = -224
$ws.Range("N97").Value = -4129.5

$ws.Range("H102").Value = 34499.03
$ws.Range("I102").Value = 1952.4117
$ws.Range("K102").Value = 1952.4117
$ws.Range("M102").Value = -330.4117000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2943634.5
$ws.Range("I136").Value = 4168540.5
$ws.Range("J136").Value = 3860
$ws.Range("K136").Value = 12505621.5
$ws.Range("L136").Value = 11580
$ws.Range("M136").Value = -12503071.5
$ws.Range("N136").Value = -16680

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws.Range("H141").Value = 30000
$ws.Range("J141").Value = 30000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = ""

$ws.Range("H62").Value = 4253.8184
$ws.Range("I62").Value = 3760
$ws.Range("J62").Value = 4665.3335
$ws.Range("K62").Value = 3760
$ws.Range("L62").Value = 4665.3335
$ws.Range("M62").Value = -3136
$ws.Range("N62").Value = -5913.3335

$ws.Range("H65").Value = 4253.8184
$ws.Range("I65").Value = 3760
$ws.Range("J65").Value = 4665.3335
$ws.Range("K65").Value = 18800
$ws.Range("L65").Value = 23326.6675
$ws.Range("M65").Value = -15680
$ws.Range("N65").Value = -29566.6675

$ws.Range("H81").Value = 1557.2142
$ws.Range("I81").Value = 1267
$ws.Range("J81").Value = 1847.4286
$ws.Range("K81").Value = 2534
$ws.Range("L81").Value = 3694.8572
$ws.Range("M81").Value = -1473
$ws.Range("N81").Value = -5816.8572

$ws.Range("H84").Value = 1557.2142
$ws.Range("I84").Value = 1267
$ws.Range("J84").Value = 1847.4286
$ws.Range("K84").Value = 12670
$ws.Range("L84").Value = 18474.286
$ws.Range("M84").Value = -7366
$ws.Range("N84").Value = -29082.286

$ws.Range("H132").Value = 226530.86
$ws.Range("I132").Value = 371599.16
$ws.Range("J132").Value = 8928.444
$ws.Range("K132").Value = 1114797.48
$ws.Range("L132").Value = 26785.332
$ws.Range("M132").Value = -1112267.48
$ws.Range("N132").Value = -31845.332

$ws.Range("H139").Value = 29700
$ws.Range("J139").Value = 29700
$ws.Range("L139").Value = 29700
$ws.Range("N139").Value = -39980

$ws.Range("H141").Value = 28850
$ws.Range("J141").Value = 28850
$ws.Range("L141").Value = 28850
$ws.Range("N141").Value = -39210

